$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions) - update "想去人数" (interest count) in column F
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 21203
$ws1.Range("F8").Value = 7909
$ws1.Range("F9").Value = 549
$ws1.Range("F12").Value = 310
$ws1.Range("F13").Value = 61
$ws1.Range("F15").Value = 165
$ws1.Range("F20").Value = 525
$ws1.Range("F21").Value = 82
$ws1.Range("F26").Value = 346
$ws1.Range("F27").Value = 1180
$ws1.Range("F28").Value = 51
$ws1.Range("F29").Value = 40
$ws1.Range("F30").Value = 221
$ws1.Range("F31").Value = 5213
$ws1.Range("F32").Value = 602
$ws1.Range("F33").Value = 7
$ws1.Range("F34").Value = 136
$ws1.Range("F35").Value = 5034
$ws1.Range("F38").Value = 43
$ws1.Range("F40").Value = 13071
$ws1.Range("F41").Value = 1362
$ws1.Range("F42").Value = 131
$ws1.Range("F43").Value = 51
$ws1.Range("F45").Value = 304
$ws1.Range("F46").Value = 429

# Sheet 2: 演出 (Performances) - update column F
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 326

# Sheet 4: 全部类型 (All types) - update column F
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 21203
$ws4.Range("F7").Value = 7909
$ws4.Range("F8").Value = 549
$ws4.Range("F11").Value = 310
$ws4.Range("F12").Value = 61
$ws4.Range("F14").Value = 165
$ws4.Range("F18").Value = 525
$ws4.Range("F19").Value = 82
$ws4.Range("F24").Value = 346
$ws4.Range("F25").Value = 1180
$ws4.Range("F26").Value = 51
$ws4.Range("F27").Value = 40
$ws4.Range("F28").Value = 221
$ws4.Range("F29").Value = 326
$ws4.Range("F30").Value = 602
$ws4.Range("F32").Value = 7
$ws4.Range("F33").Value = 136
$ws4.Range("F35").Value = 5034
$ws4.Range("F38").Value = 43
$ws4.Range("F40").Value = 13071
$ws4.Range("F41").Value = 1362
$ws4.Range("F42").Value = 131
$ws4.Range("F43").Value = 51
$ws4.Range("F45").Value = 304
$ws4.Range("F46").Value = 429
